# Add the missing day's row (row 45: 2025-12-11 / serial 46002) to the
# "Daily 100 Error Counts" sheet, and move the active selection down to it
# (matching the author's next day's tracking row: A45:D45).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A keeps the existing date-formatted style (s="2") already applied
# to A45; we only need to populate the values for the new day's counts.
$ws.Range("A45").Value = 46002
$ws.Range("B45").Value = 612
$ws.Range("C45").Value = 13
$ws.Range("D45").Value = 599

# Move the selection to the newly filled row, same shape as before (A44:D44 -> A45:D45).
$ws.Range("A45:D45").Select()
